# This script reproduces an upstream re-sync of the "Artfynd" sheet in which a
# batch of observation rows got re-keyed: the non-key cell content of several
# rows (everything except the row number itself) was rotated among a small
# set of row positions. Concretely, for each cycle of row numbers below, row
# cycle[i] ends up holding the data that used to live in row cycle[i+1]
# (wrapping around).
#
# Columns A:X and AB:AY are copied (Y/AA — the text date columns — are
# identical "2026-02-07" literals on every affected row, so they are left
# untouched rather than round-tripped through Value2, which would otherwise
# coerce that date-like text into a real date).
#
# Column I is textual on these rows (it stores things like the literal
# string "2"), so it is restored with NumberFormat "@" forced first to keep
# it from being reinterpreted as a number.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Get-RowSnapshot($row) {
    $seg1 = $ws.Range("A" + $row + ":X" + $row).Value2
    $seg2 = $ws.Range("AB" + $row + ":AY" + $row).Value2
    $iVal = $ws.Range("I" + $row).Value2
    return @($seg1, $seg2, $iVal)
}

function Set-RowSnapshot($row, $snap) {
    $ws.Range("A" + $row + ":X" + $row).Value = $snap[0]
    $ws.Range("AB" + $row + ":AY" + $row).Value = $snap[1]
    $ws.Range("I" + $row).NumberFormat = "@"
    $ws.Range("I" + $row).Value = $snap[2]
}

# Row-content rotation cycles: row cycle[i] <- old content of row cycle[i+1].
$cycles = @(
    @("3", "7", "5", "6", "4"),
    @("8", "9"),
    @("10", "13"),
    @("23", "24"),
    @("38", "39", "40", "41")
)

foreach ($cycle in $cycles) {
    # Snapshot every row in the cycle BEFORE writing any of them, since the
    # cycle is a rotation (sources and destinations overlap).
    $snapshots = @{}
    foreach ($row in $cycle) {
        $snapshots[$row] = Get-RowSnapshot $row
    }

    $n = $cycle.Length
    for ($i = 0; $i -lt $n; $i++) {
        $destRow = $cycle[$i]
        $srcRow = $cycle[($i + 1) % $n]
        Set-RowSnapshot $destRow $snapshots[$srcRow]
    }
}
